{"js": "const replacements = [\n  [\"989\u00f79=109, 8\", \"981\u00f77=140, 1\"],\n  [\"977\u00f74=244, 1\", \"135\u00f77=19, 2\"],\n  [\"721\u00f72=360, 1\", \"387\u00f78=48, 3\"],\n  [\"809\u00f79=89, 8\", \"155\u00f73=51, 2\"],\n  [\"599\u00f72=299, 1\", \"600\u00f72=300, 0\"],\n  [\"888\u00f78=111, 0\", \"261\u00f72=130, 1\"],\n  [\"431\u00f73=143, 2\", \"535\u00f73=178, 1\"],\n  [\"727\u00f73=242, 1\", \"869\u00f77=124, 1\"],\n  [\"762\u00f76=127, 0\", \"249\u00f76=41, 3\"],\n  [\"338\u00f78=42, 2\", \"266\u00f78=33, 2\"],\n  [\"515\u00f74=128, 3\", \"775\u00f73=258, 1\"],\n  [\"553\u00f77=79, 0\", \"555\u00f77=79, 2\"],\n  [\"192\u00f74=48, 0\", \"789\u00f75=157, 4\"],\n  [\"532\u00f76=88, 4\", \"931\u00f78=116, 3\"],\n  [\"613\u00f77=87, 4\", \"503\u00f74=125, 3\"],\n  [\"352\u00f73=117, 1\", \"419\u00f76=69, 5\"],\n  [\"919\u00f74=229, 3\", \"136\u00f76=22, 4\"],\n  [\"158\u00f79=17, 5\", \"995\u00f76=165, 5\"],\n  [\"959\u00f72=479, 1\", \"773\u00f77=110, 3\"],\n  [\"996\u00f73=332, 0\", \"130\u00f74=32, 2\"],\n  [\"927\u00f76=154, 3\", \"429\u00f72=214, 1\"],\n  [\"410\u00f79=45, 5\", \"889\u00f75=177, 4\"],\n  [\"262\u00f77=37, 3\", \"253\u00f78=31, 5\"],\n  [\"528\u00f76=88, 0\", \"740\u00f77=105, 5\"],\n  [\"202\u00f72=101, 0\", \"642\u00f73=214, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"989\u00f79=109, 8\", \"981\u00f77=140, 1\"),\n    @(\"977\u00f74=244, 1\", \"135\u00f77=19, 2\"),\n    @(\"721\u00f72=360, 1\", \"387\u00f78=48, 3\"),\n    @(\"809\u00f79=89, 8\", \"155\u00f73=51, 2\"),\n    @(\"599\u00f72=299, 1\", \"600\u00f72=300, 0\"),\n    @(\"888\u00f78=111, 0\", \"261\u00f72=130, 1\"),\n    @(\"431\u00f73=143, 2\", \"535\u00f73=178, 1\"),\n    @(\"727\u00f73=242, 1\", \"869\u00f77=124, 1\"),\n    @(\"762\u00f76=127, 0\", \"249\u00f76=41, 3\"),\n    @(\"338\u00f78=42, 2\", \"266\u00f78=33, 2\"),\n    @(\"515\u00f74=128, 3\", \"775\u00f73=258, 1\"),\n    @(\"553\u00f77=79, 0\", \"555\u00f77=79, 2\"),\n    @(\"192\u00f74=48, 0\", \"789\u00f75=157, 4\"),\n    @(\"532\u00f76=88, 4\", \"931\u00f78=116, 3\"),\n    @(\"613\u00f77=87, 4\", \"503\u00f74=125, 3\"),\n    @(\"352\u00f73=117, 1\", \"419\u00f76=69, 5\"),\n    @(\"919\u00f74=229, 3\", \"136\u00f76=22, 4\"),\n    @(\"158\u00f79=17, 5\", \"995\u00f76=165, 5\"),\n    @(\"959\u00f72=479, 1\", \"773\u00f77=110, 3\"),\n    @(\"996\u00f73=332, 0\", \"130\u00f74=32, 2\"),\n    @(\"927\u00f76=154, 3\", \"429\u00f72=214, 1\"),\n    @(\"410\u00f79=45, 5\", \"889\u00f75=177, 4\"),\n    @(\"262\u00f77=37, 3\", \"253\u00f78=31, 5\"),\n    @(\"528\u00f76=88, 0\", \"740\u00f77=105, 5\"),\n    @(\"202\u00f72=101, 0\", \"642\u00f73=214, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
